# NIT-9013909029.xlsx - "Actualiza base de datos EC y agrega parte 1 de
# nuevos estado de cuenta"
#
# The worker detail table (rows 16-30, columns B:G) is re-sorted/updated:
#  - Periods are regrouped per worker (LUIS JAVIER ARRIETA YEPEZ /
#    WILBER RAFAEL ARRIETA YEPEZ), now interleaved and in ascending period
#    order (2102, 2106, 2108, 2109, 2110, 2111, 2112, 2201).
#  - WILBER's "Salario Basico" (col G) changes from 908526 to 877803.
#  - A couple of "Valor Mora" (col F) values change for period 2102/2110.
#
# Columns: B=Tipo Doc, C=N Doc Trabajador, D=Nombre Trabajador,
#          E=Periodo Mora, F=Valor Mora, G=Salario Basico

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 16; B = "CC"; C = "1051824044"; D = "LUIS JAVIER ARRIETA YEPEZ";   E = "2102"; F = 15748; G = 908526 },
    @{ Row = 17; B = "CC"; C = "1051824044"; D = "LUIS JAVIER ARRIETA YEPEZ";   E = "2106"; F = 36341; G = 908526 },
    @{ Row = 18; B = "CC"; C = "1143348214"; D = "WILBER RAFAEL ARRIETA YEPEZ"; E = "2106"; F = 36341; G = 877803 },
    @{ Row = 19; B = "CC"; C = "1051824044"; D = "LUIS JAVIER ARRIETA YEPEZ";   E = "2108"; F = 36341; G = 908526 },
    @{ Row = 20; B = "CC"; C = "1143348214"; D = "WILBER RAFAEL ARRIETA YEPEZ"; E = "2108"; F = 36341; G = 877803 },
    @{ Row = 21; B = "CC"; C = "1051824044"; D = "LUIS JAVIER ARRIETA YEPEZ";   E = "2109"; F = 36341; G = 908526 },
    @{ Row = 22; B = "CC"; C = "1143348214"; D = "WILBER RAFAEL ARRIETA YEPEZ"; E = "2109"; F = 36341; G = 877803 },
    @{ Row = 23; B = "CC"; C = "1051824044"; D = "LUIS JAVIER ARRIETA YEPEZ";   E = "2110"; F = 36341; G = 908526 },
    @{ Row = 24; B = "CC"; C = "1143348214"; D = "WILBER RAFAEL ARRIETA YEPEZ"; E = "2110"; F = 36341; G = 877803 },
    @{ Row = 25; B = "CC"; C = "1051824044"; D = "LUIS JAVIER ARRIETA YEPEZ";   E = "2111"; F = 36341; G = 908526 },
    @{ Row = 26; B = "CC"; C = "1143348214"; D = "WILBER RAFAEL ARRIETA YEPEZ"; E = "2111"; F = 36341; G = 877803 },
    @{ Row = 27; B = "CC"; C = "1051824044"; D = "LUIS JAVIER ARRIETA YEPEZ";   E = "2112"; F = 36341; G = 908526 },
    @{ Row = 28; B = "CC"; C = "1143348214"; D = "WILBER RAFAEL ARRIETA YEPEZ"; E = "2112"; F = 36341; G = 877803 },
    @{ Row = 29; B = "CC"; C = "1051824044"; D = "LUIS JAVIER ARRIETA YEPEZ";   E = "2201"; F = 30284; G = 908526 },
    @{ Row = 30; B = "CC"; C = "1143348214"; D = "WILBER RAFAEL ARRIETA YEPEZ"; E = "2201"; F = 29260; G = 877803 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
}
